$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the changed cells keep their original text (string) representation
# rather than being auto-converted to numbers by Excel, by forcing a Text
# number format on each target cell before assigning its new value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.113.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.483.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.05"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.44%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.483.26"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.923.54"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.45"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.031.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.478.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.32%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.27"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.27"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.02"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.28"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.69"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.23"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.81"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.327"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.77"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.86"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.514"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.61%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.13%  "
